$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SCD0338")

# Update the test marker / notes / date cells in row 3
$ws.Range("L3").Value = "Test Marker"
$ws.Range("M3").Value = "notes"
$ws.Range("P3").Value = "2022-10-28"

# Reset the view: scroll back so column A is visible again and select A3
$ws.Activate()
$ws.Range("A3").Select()
